$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.835.11"

# Row 3
$ws.Range("D3").Value = "2.565.70"
$ws.Range("E3").Value = "  +1.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "'310.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6
$ws.Range("D6").Value = "'98.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.13%  "

# Row 7
$ws.Range("D7").Value = "'0.572"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  +0.08%  "

# Row 10
$ws.Range("D10").Value = "'35.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "

# Row 11
$ws.Range("D11").Value = "'0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "

# Row 12
$ws.Range("D12").Value = "'7.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13
$ws.Range("D13").Value = "2.959.60"
$ws.Range("E13").Value = "  +1.30%  "

# Row 14
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("D15").Value = "'15.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.29%  "

# Row 16
$ws.Range("D16").Value = "2.474.43"
$ws.Range("E16").Value = "  -2.81%  "

# Row 17
$ws.Range("E17").Value = "  -0.87%  "

# Row 18
$ws.Range("D18").Value = "42.848.50"
$ws.Range("E18").Value = "  +0.15%  "

# Row 19
$ws.Range("E19").Value = "  -1.19%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.27%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("D22").Value = "'69.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "

# Row 23
$ws.Range("D23").Value = "'248.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "

# Row 24
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("D26").Value = "'26.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.14%  "

# Row 27
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("D28").Value = "'2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "

# Row 29
$ws.Range("D29").Value = "'39.96"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'10.18"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'158.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.72%  "

# Row 32
$ws.Range("D32").Value = "'5.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.37%  "

# Row 33
$ws.Range("D33").Value = "'0.0798"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.54%  "

# Row 34
$ws.Range("D34").Value = "'2.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.10%  "

# Row 35
$ws.Range("E35").Value = "  -1.42%  "

# Row 36
$ws.Range("E36").Value = "  -1.03%  "

# Row 37
$ws.Range("D37").Value = "'18.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.00%  "

# Row 38
$ws.Range("E38").Value = "  +13.01%  "

# Row 39
$ws.Range("E39").Value = "  +0.78%  "

# Row 40
$ws.Range("E40").Value = "  -0.20%  "

# Row 41
$ws.Range("D41").Value = "'22.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.09%  "

# Row 42
$ws.Range("E42").Value = "  +7.24%  "

# Row 43
$ws.Range("E43").Value = "  -0.08%  "

# Row 44
$ws.Range("E44").Value = "  -0.66%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.993.98"
$ws.Range("E45").Value = "  -1.78%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "

# Row 47
$ws.Range("D47").Value = "'9.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.47%  "

# Row 48
$ws.Range("D48").Value = "2.811.67"
$ws.Range("E48").Value = "  +1.34%  "

# Row 49
$ws.Range("D49").Value = "'0.194"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.71%  "

# Row 50
$ws.Range("D50").Value = "'81.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.59%  "

# Row 51
$ws.Range("D51").Value = "'73.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
